$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally holds a one-row header (A1:J1) plus one data row
# (A2:J2). Katalon AI regenerated the sheet down to a single blank,
# still-formatted placeholder cell at A1 - everything else goes away.

# Drop the whole data row (row 2).
$ws.Rows("2:2").Clear()

# Drop the rest of the header row, keeping only column A.
$ws.Range("B1:J1").Clear()

# A1 keeps its header style, but loses its text - it becomes an empty cell.
$ws.Range("A1").ClearContents()
